$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new venue/county rows below the existing data (rows 128-130)
$ws.Range("A128").Value = "Páirc Grattan, Inniskeen"
$ws.Range("B128").Value = "Monaghan"

$ws.Range("A129").Value = "Toomevara, Tipperary"
$ws.Range("B129").Value = "Tipperary"

$ws.Range("A130").Value = "Corrigan Park"
$ws.Range("B130").Value = "Antrim"

# Reflect the new scroll position / selection that was saved with the sheet
$excel.ActiveWindow.ScrollRow = 118
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B131").Select()
